$wb = $excel.ActiveWorkbook

# Sheet order: 1 = "Semilla 6", 2 = "Semilla 3"
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet "Semilla 6": add a new row 14 (previously hard-coded portin data) ---
$ws1.Cells.Item(14, 1).NumberFormat = "@"
$ws1.Cells.Item(14, 1).Value = "10960370"

$ws1.Cells.Item(14, 2).NumberFormat = "@"
$ws1.Cells.Item(14, 2).Value = "884243417"

$ws1.Cells.Item(14, 3).NumberFormat = "@"
$ws1.Cells.Item(14, 3).Value = "3043209868"

$ws1.Cells.Item(14, 4).NumberFormat = "@"
$ws1.Cells.Item(14, 4).Value = "732111324707278 "

# --- Sheet "Semilla 3": adjust row 13 and add new row 14 ---
$ws2.Cells.Item(13, 4).NumberFormat = "@"
$ws2.Cells.Item(13, 4).Value = "732111324707278"

$ws2.Cells.Item(14, 1).NumberFormat = "@"
$ws2.Cells.Item(14, 1).Value = "10960370"

$ws2.Cells.Item(14, 2).NumberFormat = "@"
$ws2.Cells.Item(14, 2).Value = "884243417"

$ws2.Cells.Item(14, 3).NumberFormat = "@"
$ws2.Cells.Item(14, 3).Value = "3043209863"

$ws2.Cells.Item(14, 4).NumberFormat = "@"
$ws2.Cells.Item(14, 4).Value = "732111324707277"

# --- Update selections / active sheet ---
$ws1.Range("A13:B14").Select()
$ws2.Range("C15").Select()
